# Added test suite file and updated test method
# - Fill in two new values on the "Invalid" sheet (A3="ghvg", B2="vjhv"),
#   which also appends two new shared strings.
# - Make "Invalid" the active/selected sheet (tab) instead of "Valid",
#   with its selection on B2.

$wb = $excel.ActiveWorkbook
$wsInvalid = $wb.Worksheets.Item("Invalid")

# Order matters for shared-string insertion order: A3 ("ghvg") must be
# written before B2 ("vjhv") so the new <si> entries land at indices 8/9
# respectively, matching the target workbook.
$wsInvalid.Range("A3").Value = "ghvg"
$wsInvalid.Range("B2").Value = "vjhv"

# Switch the active tab to "Invalid" and select B2 there (this clears
# tabSelected from the previously active "Valid" sheet).
$wsInvalid.Select()
$wsInvalid.Range("B2").Select()
